$wb = $excel.ActiveWorkbook

# Delete column C ("VevoLab [mm^3]") on sheets 1339, 1340, 1341, 1342,
# shifting D,E,F,G left to C,D,E,F.
foreach ($name in @("1339", "1340", "1341", "1342")) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Columns.Item(3).Delete()
}

# Sheet 1341 ("1341") had a stray applied-number-format on what was D2;
# after the shift that cell (now C2) loses that formatting entirely.
$ws4 = $wb.Worksheets.Item("1341")
$ws4.Range("C2").ClearFormats()

# Restore per-sheet selections to what they were left at after editing.
$ws2 = $wb.Worksheets.Item("1339")
$ws2.Range("C1:C1048576").Select()

$ws3 = $wb.Worksheets.Item("1340")
$ws3.Range("C1:C1048576").Select()

$ws4.Range("C1:C1048576").Select()

$ws5 = $wb.Worksheets.Item("1342")
$ws5.Activate()
$ws5.Range("F14").Select()
